# Adapt column header formatting to respective input file names.
# "*_old" headers become "*_FV2404", "*_new" headers become "*_FV2410",
# wrap the data range in an Excel Table ("Table1") and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @(
  "Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404",
  "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404",
  "diff",
  "Segmentname_FV2410", "Segmentgruppe_FV2410", "Segment_FV2410", "Datenelement_FV2410", "Segment ID_FV2410",
  "Code_FV2410", "Qualifier_FV2410", "Beschreibung_FV2410", "Bedingungsausdruck_FV2410", "Bedingung_FV2410"
)

# Rename the header row (row 1) in place so the shared-string entries are
# updated rather than duplicated.
for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Wrap A1:U84 in a real Excel table named "Table1" with an autofilter,
# matching the worksheet's used range.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U84"), $null, 1)
$tbl.Name = "Table1"

# Freeze the header row: select A2 and freeze panes above/left of it.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
